$d = $word.ActiveDocument

# --- 1) "Known Issues" heading: merge the split runs (and drop the
#        gramStart/gramEnd proofErr markers left over from a grammar
#        check) into a single clean run. ---
$found = $false
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "Known Issues") {
        $found = $true
        # Insert a fresh, clean paragraph (inherits the Heading1 style)
        # immediately before the old one, then delete the old one. This
        # removes the stray w:proofErr elements and the split runs.
        # Note: after InsertBefore, $p.Range itself now refers to the
        # *newly inserted* paragraph, so the original (messy) paragraph
        # is the next one.
        $p.Range.InsertBefore("Known Issues`r")
        $old = $p.Next()
        $old.Range.Delete()
        break
    }
}

# --- 2) Add a new sub-bullet "Partial load of queue and data" right
#        after the "Dealing with images better..." bullet. ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Dealing with images better")) {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Partial load of queue and data"
        break
    }
}
